# Renames the boolean "centralized" (NO/YES) column to a "scale"
# (NONE/BUILDING/CITY/DISTRICT) column across the DHW, HEATING, COOLING
# and ELECTRICITY sheets, per discussion feedback on the LCA database.

$wb = $excel.ActiveWorkbook

# --- DHW ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("DHW")
$ws.Range("D1").Value = "scale"
$ws.Range("D2").Value = "NONE"
$ws.Range("D3").Value = "BUILDING"
$ws.Range("D4").Value = "BUILDING"
$ws.Range("A16").Select() | Out-Null

# --- HEATING -------------------------------------------------------------
$ws = $wb.Worksheets.Item("HEATING")
$ws.Range("D1").Value = "scale"
$ws.Range("D2").Value = "NONE"
$ws.Range("D3").Value = "BUILDING"
$ws.Range("D4").Value = "BUILDING"
$ws.Range("C8").Select() | Out-Null

# --- ELECTRICITY -----------------------------------------------------------
$ws = $wb.Worksheets.Item("ELECTRICITY")
$ws.Range("D1").Value = "scale"
$ws.Range("D2").Value = "NONE"
$ws.Range("D3").Value = "BUILDING"
$ws.Range("D4").Value = "CITY"
$ws.Range("D4").Select() | Out-Null

# --- COOLING -------------------------------------------------------------
$ws = $wb.Worksheets.Item("COOLING")
$ws.Range("D1").Value = "scale"
$ws.Range("D2").Value = "NONE"
$ws.Range("D3").Value = "BUILDING"
$ws.Range("D4").Value = "BUILDING"
$ws.Range("D5").Value = "DISTRICT"
$ws.Range("D6").Value = "DISTRICT"
$ws.Range("D7").Value = "DISTRICT"
$ws.Range("C11").Select() | Out-Null

# --- FUELS -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("FUELS")
$ws.Range("C16").Select() | Out-Null

# ELECTRICITY becomes the active/selected tab (was FUELS).
$wb.Worksheets.Item("ELECTRICITY").Activate() | Out-Null
